$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 12.36652176781263
$ws.Range("C2").Value = 6.710362301308201
$ws.Range("E2").Value = 11.74266040059292
$ws.Range("F2").Value = 16.86991607391233
$ws.Range("G2").Value = 44.99892709505789
$ws.Range("H2").Value = 18.45801825050868
$ws.Range("I2").Value = 29.30628798620206
$ws.Range("K2").Value = 10.08830146424279
$ws.Range("L2").Value = 10.31677121513143
$ws.Range("M2").Value = 14.3987000792445
$ws.Range("B3").Value = 12.18046356039775
$ws.Range("C3").Value = 6.63703623605658
$ws.Range("E3").Value = 11.76408651046654
$ws.Range("F3").Value = 15.89584955866808
$ws.Range("G3").Value = 45.13230950563796
$ws.Range("H3").Value = 18.51750605185546
$ws.Range("I3").Value = 29.41173216214271
$ws.Range("K3").Value = 9.959583861590385
$ws.Range("L3").Value = 10.32630042546873
$ws.Range("M3").Value = 14.37802520705135
$ws.Range("B4").Value = 12.06729538813024
$ws.Range("C4").Value = 6.590801366029492
$ws.Range("E4").Value = 11.77860313348295
$ws.Range("F4").Value = 15.26997757108491
$ws.Range("G4").Value = 45.22582510834144
$ws.Range("H4").Value = 18.55689640392668
$ws.Range("I4").Value = 29.48155766062121
$ws.Range("K4").Value = 9.88167543374064
$ws.Range("L4").Value = 10.33356965493217
$ws.Range("M4").Value = 14.36749306120887
$ws.Range("B5").Value = 12.02150646181504
$ws.Range("C5").Value = 6.571662680592414
$ws.Range("E5").Value = 11.78486134521437
$ws.Range("F5").Value = 15.00819731993403
$ws.Range("G5").Value = 45.26684364524506
$ws.Range("H5").Value = 18.57366833102252
$ws.Range("I5").Value = 29.51128892823386
$ws.Range("K5").Value = 9.850246678531182
$ws.Range("L5").Value = 10.33688896805814
$ws.Range("M5").Value = 14.36374849369082
$ws.Range("B6").Value = 12.01392475987686
$ws.Range("C6").Value = 6.56846695728732
$ws.Range("E6").Value = 11.78592121717912
$ws.Range("F6").Value = 14.96433081551593
$ws.Range("G6").Value = 45.27383016703286
$ws.Range("H6").Value = 18.57649677700219
$ws.Range("I6").Value = 29.51630286429825
$ws.Range("K6").Value = 9.845048322700718
$ws.Range("L6").Value = 10.33746171278329
$ws.Range("M6").Value = 14.36315987072918
$ws.Range("B7").Value = 12.06667645749996
$ws.Range("C7").Value = 6.590544450285933
$ws.Range("E7").Value = 11.77868614638965
$ws.Range("F7").Value = 15.26647399323137
$ws.Range("G7").Value = 45.22636653095607
$ws.Range("H7").Value = 18.55711968098711
$ws.Range("I7").Value = 29.48195345872662
$ws.Range("K7").Value = 9.881250232478207
$ws.Range("L7").Value = 10.33361297419806
$ws.Range("M7").Value = 14.36744033982882
$ws.Range("B8").Value = 12.3021813797321
$ws.Range("C8").Value = 6.685335400644325
$ws.Range("E8").Value = 11.74976588683975
$ws.Range("F8").Value = 16.53996406344768
$ws.Range("G8").Value = 45.04249932470844
$ws.Range("H8").Value = 18.47793492392567
$ws.Range("I8").Value = 29.34158958388959
$ws.Range("K8").Value = 10.04370837791688
$ws.Range("L8").Value = 10.31976279145692
$ws.Range("M8").Value = 14.39112486973961
$ws.Range("B9").Value = 12.76995486086559
$ws.Range("C9").Value = 6.861292757831203
$ws.Range("E9").Value = 11.70383550535475
$ws.Range("F9").Value = 19.00274580682531
$ws.Range("G9").Value = 44.774606168392
$ws.Range("H9").Value = 18.34539600779572
$ws.Range("I9").Value = 29.10671554591397
$ws.Range("K9").Value = 10.36958212796989
$ws.Range("L9").Value = 10.30383619807299
$ws.Range("M9").Value = 14.45455519761272
$ws.Range("B10").Value = 13.11386782526055
$ws.Range("C10").Value = 6.984121051799814
$ws.Range("E10").Value = 11.67664166546587
$ws.Range("F10").Value = 20.67494806633232
$ws.Range("G10").Value = 44.63490641835987
$ws.Range("H10").Value = 18.26190133287403
$ws.Range("I10").Value = 28.95884203137886
$ws.Range("K10").Value = 10.61126197182876
$ws.Range("L10").Value = 10.29895441166589
$ws.Range("M10").Value = 14.51126400120465
$ws.Range("B11").Value = 13.26967883072583
$ws.Range("C11").Value = 7.038519045303564
$ws.Range("E11").Value = 11.66568811395777
$ws.Range("F11").Value = 21.3917225636224
$ws.Range("G11").Value = 44.58387769114042
$ws.Range("H11").Value = 18.22693564604424
$ws.Range("I11").Value = 28.8969472792312
$ws.Range("K11").Value = 10.72124104691027
$ws.Range("L11").Value = 10.29820686992018
$ws.Range("M11").Value = 14.53919755484465
$ws.Range("B12").Value = 13.32853297119597
$ws.Range("C12").Value = 7.058899227859026
$ws.Range("E12").Value = 11.66174363399106
$ws.Range("F12").Value = 21.65686569030329
$ws.Range("G12").Value = 44.56636374371083
$ws.Range("H12").Value = 18.21412920711033
$ws.Range("I12").Value = 28.87428350813267
$ws.Range("K12").Value = 10.76285524924858
$ws.Range("L12").Value = 10.2981348662045
$ws.Range("M12").Value = 14.55007692541036
$ws.Range("B13").Value = 13.31586514422911
$ws.Range("C13").Value = 7.054519842660191
$ws.Range("E13").Value = 11.6625841086359
$ws.Range("F13").Value = 21.60004134736742
$ws.Range("G13").Value = 44.57005507223386
$ws.Range("H13").Value = 18.21686797820021
$ws.Range("I13").Value = 28.879130086609
$ws.Range("K13").Value = 10.75389492443806
$ws.Range("L13").Value = 10.29814099900733
$ws.Range("M13").Value = 14.54772054079969
$ws.Range("B14").Value = 13.27452408104784
$ws.Range("C14").Value = 7.040200165148254
$ws.Range("E14").Value = 11.66535952516938
$ws.Range("F14").Value = 21.4136618050453
$ws.Range("G14").Value = 44.58240049404401
$ws.Range("H14").Value = 18.22587334529862
$ws.Range("I14").Value = 28.89506718665683
$ws.Range("K14").Value = 10.72466553720569
$ws.Range("L14").Value = 10.2981967205795
$ws.Range("M14").Value = 14.54008659747754
$ws.Range("B15").Value = 13.24918055548512
$ws.Range("C15").Value = 7.031400212674975
$ws.Range("E15").Value = 11.66708602404692
$ws.Range("F15").Value = 21.29868154950795
$ws.Range("G15").Value = 44.59019832743512
$ws.Range("H15").Value = 18.2314459685924
$ws.Range("I15").Value = 28.90493002474884
$ws.Range("K15").Value = 10.70675635391301
$ws.Range("L15").Value = 10.29825831554303
$ws.Range("M15").Value = 14.5354496825817
$ws.Range("B16").Value = 13.10366751702623
$ws.Range("C16").Value = 6.980535708485509
$ws.Range("E16").Value = 11.67738598868928
$ws.Range("F16").Value = 20.62722412089977
$ws.Range("G16").Value = 44.6384940121101
$ws.Range("H16").Value = 18.26424716808432
$ws.Range("I16").Value = 28.96299527601665
$ws.Range("K16").Value = 10.60407207781148
$ws.Range("L16").Value = 10.29903284577875
$ws.Range("M16").Value = 14.50948098995036
$ws.Range("B17").Value = 13.01419603261933
$ws.Range("C17").Value = 6.948948765460456
$ws.Range("E17").Value = 11.68406736253517
$ws.Range("F17").Value = 20.20408069597325
$ws.Range("G17").Value = 44.67133550283195
$ws.Range("H17").Value = 18.28514261548449
$ws.Range("I17").Value = 28.99999395711596
$ws.Range("K17").Value = 10.54106100028162
$ws.Range("L17").Value = 10.29988479407011
$ws.Range("M17").Value = 14.49409328774633
$ws.Range("B18").Value = 12.96267773200013
$ws.Range("C18").Value = 6.93064213706554
$ws.Range("E18").Value = 11.68804371650597
$ws.Range("F18").Value = 19.95656407809801
$ws.Range("G18").Value = 44.6914031442357
$ws.Range("H18").Value = 18.29744501471772
$ws.Range("I18").Value = 29.02178031076386
$ws.Range("K18").Value = 10.50482445982393
$ws.Range("L18").Value = 10.30051352838695
$ws.Range("M18").Value = 14.48544413893546
$ws.Range("B19").Value = 12.9452265435869
$ws.Range("C19").Value = 6.924420222347138
$ws.Range("E19").Value = 11.68941296525843
$ws.Range("F19").Value = 19.87204792380568
$ws.Range("G19").Value = 44.69839975417989
$ws.Range("H19").Value = 18.30165913768493
$ws.Range("I19").Value = 29.02924359694616
$ws.Range("K19").Value = 10.49255754988925
$ws.Range("L19").Value = 10.30075025392071
$ws.Range("M19").Value = 14.48255045308556
$ws.Range("B20").Value = 13.02372670600361
$ws.Range("C20").Value = 6.952325649214219
$ws.Range("E20").Value = 11.68334231472217
$ws.Range("F20").Value = 20.24955283636154
$ws.Range("G20").Value = 44.66771747528328
$ws.Range("H20").Value = 18.28288887387768
$ws.Range("I20").Value = 28.99600303490412
$ws.Range("K20").Value = 10.54776831795614
$ws.Range("L20").Value = 10.29977975146157
$ws.Range("M20").Value = 14.49571052570331
$ws.Range("B21").Value = 13.28667141331724
$ws.Range("C21").Value = 7.044412201895647
$ws.Range("E21").Value = 11.66453880116153
$ws.Range("F21").Value = 21.46857628470577
$ws.Range("G21").Value = 44.57872516690086
$ws.Range("H21").Value = 18.22321645933983
$ws.Range("I21").Value = 28.89036504266634
$ws.Range("K21").Value = 10.73325209582425
$ws.Range("L21").Value = 10.29817463183673
$ws.Range("M21").Value = 14.54232073469646
$ws.Range("B22").Value = 13.45763298981475
$ws.Range("C22").Value = 7.103315719226851
$ws.Range("E22").Value = 11.65343495521162
$ws.Range("F22").Value = 22.22866616901552
$ws.Range("G22").Value = 44.53111443597626
$ws.Range("H22").Value = 18.18674882036906
$ws.Range("I22").Value = 28.82583926529598
$ws.Range("K22").Value = 10.85427023584828
$ws.Range("L22").Value = 10.29835548321108
$ws.Range("M22").Value = 14.57453796462816
$ws.Range("B23").Value = 13.36648696386902
$ws.Range("C23").Value = 7.071997120082746
$ws.Range("E23").Value = 11.65925295830886
$ws.Range("F23").Value = 21.82633154458858
$ws.Range("G23").Value = 44.55555701589132
$ws.Range("H23").Value = 18.20598045100356
$ws.Range("I23").Value = 28.8598642212712
$ws.Range("K23").Value = 10.78971171636097
$ws.Range("L23").Value = 10.29814669893381
$ws.Range("M23").Value = 14.55718441389292
$ws.Range("B24").Value = 13.01941813326471
$ws.Range("C24").Value = 6.95079941707248
$ws.Range("E24").Value = 11.68366968786009
$ws.Range("F24").Value = 20.22900810905287
$ws.Range("G24").Value = 44.6693494898764
$ws.Range("H24").Value = 18.28390688922769
$ws.Range("I24").Value = 28.99780572573282
$ws.Range("K24").Value = 10.54473596984442
$ws.Range("L24").Value = 10.29982680839089
$ws.Range("M24").Value = 14.49497875696668
$ws.Range("B25").Value = 12.64312633273737
$ws.Range("C25").Value = 6.814796021334754
$ws.Range("E25").Value = 11.71510878719103
$ws.Range("F25").Value = 18.34778573295695
$ws.Range("G25").Value = 44.83709007253414
$ws.Range("H25").Value = 18.37881481547063
$ws.Range("I25").Value = 29.16592492885124
$ws.Range("K25").Value = 10.28087027872123
$ws.Range("L25").Value = 10.30694482547639
$ws.Range("M25").Value = 14.43560156337159
